$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression - only B2 changes slightly
$ws.Range("B2").Value = 0.8860399364920633

# Row 3: RandomForestRegressor - B3, C3, D3 change
$ws.Range("B3").Value = 0.9948158510317802
$ws.Range("C3").Value = 0.994850533403372
$ws.Range("D3").Value = 0.9550569337232394

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, values change
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9937780421949309
$ws.Range("C4").Value = 0.9938701477839333
$ws.Range("D4").Value = 0.873371552548846

# Row 5: AdaBoostRegressor -> MLPRegressor, values change
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9965814453913229
$ws.Range("C5").Value = 0.9965925276141196
$ws.Range("D5").Value = 0.9950614347314789
